{"js": "// Edit 1: In the \"filtrado por fecha\" feedback paragraph, insert a new\n// sentence about requesting icons/images right after \"...todo mezclado.\"\n// and before \" Se dej\u00f3 2 d\u00edas...\".\n{\n  const results = context.document.body.search(\"todo mezclado.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \" Adem\u00e1s se solicit\u00f3 que se colocaran iconos e im\u00e1genes para dar una idea m\u00e1s aproximada de como seria finalmente la aplicaci\u00f3n.\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n  }\n}\n\n// Edit 2: In the \"autenticaci\u00f3n por Facebook\" feedback paragraph, insert a\n// new sentence about adding images/icons right after \"...por Facebook.\"\n// and before \" 2 d\u00edas despu\u00e9s...\".\n{\n  const results = context.document.body.search(\"autenticaci\u00f3n por Facebook.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \" Adem\u00e1s, se indic\u00f3 que deber\u00edan agregarse las im\u00e1genes e iconos de la aplicaci\u00f3n para que el cliente se haga una idea de c\u00f3mo va a quedar.\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n  }\n}\n\n// Edit 3: Remove the standalone manual-page-break paragraph that sits\n// between \"...para otro sprint.\" and the \"Sprint 2\" heading.\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  let targetIndex = -1;\n  for (let i = 0; i < paras.items.length - 1; i++) {\n    if (paras.items[i].text === \"\\f\" && paras.items[i + 1].text === \"Sprint 2\") {\n      targetIndex = i;\n      break;\n    }\n  }\n\n  if (targetIndex !== -1) {\n    paras.items[targetIndex].delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Edit 1: In the \"filtrado por fecha\" feedback paragraph, insert a new\n# sentence about requesting icons/images right after \"...todo mezclado.\"\n# and before \" Se dej\u00f3 2 d\u00edas...\".\n$find1 = $d.Content\n$find1.Find.Text = \"todo mezclado.\"\n$find1.Find.Execute() | Out-Null\nif ($find1.Find.Found) {\n    $find1.InsertAfter(\" Adem\u00e1s se solicit\u00f3 que se colocaran iconos e im\u00e1genes para dar una idea m\u00e1s aproximada de como seria finalmente la aplicaci\u00f3n.\")\n}\n\n# Edit 2: In the \"autenticaci\u00f3n por Facebook\" feedback paragraph, insert a\n# new sentence about adding images/icons right after \"...por Facebook.\"\n# and before \" 2 d\u00edas despu\u00e9s...\".\n$find2 = $d.Content\n$find2.Find.Text = \"autenticaci\u00f3n por Facebook.\"\n$find2.Find.Execute() | Out-Null\nif ($find2.Find.Found) {\n    $find2.InsertAfter(\" Adem\u00e1s, se indic\u00f3 que deber\u00edan agregarse las im\u00e1genes e iconos de la aplicaci\u00f3n para que el cliente se haga una idea de c\u00f3mo va a quedar.\")\n}\n\n# Edit 3: Remove the standalone manual-page-break paragraph that sits\n# between \"...para otro sprint.\" and the \"Sprint 2\" heading.\n$paras = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -lt $paras.Count; $i++) {\n    $cur = $paras.Item($i).Range.Text\n    $stripped = $cur -replace \"[\\f\\r\\a\\v]\", \"\"\n    if ($stripped -eq \"\" -and $cur.Length -gt 0) {\n        $nextText = $paras.Item($i + 1).Range.Text\n        if ($nextText -match \"^Sprint 2\") {\n            $targetIndex = $i\n            break\n        }\n    }\n}\nif ($targetIndex -ne -1) {\n    $paras.Item($targetIndex).Range.Delete()\n}\n\nWrite-Output \"done\"\n"}
